# Auto-generated edit script: updates cached market-value columns (H-N)
# across multiple Leve-profit sheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 549.6667
$ws.Range("I9").Value = 398
$ws.Range("K9").Value = 398
$ws.Range("M9").Value = -229
$ws.Range("H33").Value = 887.94446
$ws.Range("I33").Value = 429
$ws.Range("K33").Value = 429
$ws.Range("M33").Value = -200
$ws.Range("H64").Value = 21882.166
$ws.Range("I64").Value = 25843
$ws.Range("K64").Value = 25843
$ws.Range("M64").Value = -25595
$ws.Range("H67").Value = 21882.166
$ws.Range("I67").Value = 25843
$ws.Range("K67").Value = 25843
$ws.Range("M67").Value = -24985
$ws.Range("I106").Value = 3489
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3489
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2858
$ws.Range("N106").ClearContents()
$ws.Range("H111").Value = 656.8889
$ws.Range("I111").Value = 627.75
$ws.Range("J111").Value = 890
$ws.Range("K111").Value = 1883.25
$ws.Range("L111").Value = 2670
$ws.Range("M111").Value = 1183.75
$ws.Range("N111").Value = -8804
$ws.Range("H132").Value = 4632142.5
$ws.Range("I132").Value = 5743626.5
$ws.Range("J132").Value = 957.5
$ws.Range("K132").Value = 17230879.5
$ws.Range("L132").Value = 2872.5
$ws.Range("M132").Value = -17228349.5
$ws.Range("N132").Value = -7932.5
$ws.Range("H135").Value = 1703.0625
$ws.Range("I135").Value = 711.53845
$ws.Range("K135").Value = 6403.84605
$ws.Range("M135").Value = -3868.84605

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27830.2
$ws.Range("I32").Value = 29058.105
$ws.Range("K32").Value = 29058.105
$ws.Range("M32").Value = -28771.105
$ws.Range("H74").Value = 382988.12
$ws.Range("I74").Value = 546637.4
$ws.Range("J74").Value = 22959.8
$ws.Range("K74").Value = 546637.4
$ws.Range("L74").Value = 22959.8
$ws.Range("M74").Value = -545763.4
$ws.Range("N74").Value = -24707.8
$ws.Range("H77").Value = 382988.12
$ws.Range("I77").Value = 546637.4
$ws.Range("J77").Value = 22959.8
$ws.Range("K77").Value = 2733187
$ws.Range("L77").Value = 114799
$ws.Range("M77").Value = -2728819
$ws.Range("N77").Value = -123535
$ws.Range("H122").Value = 2776.5625
$ws.Range("I122").Value = 2897.5454
$ws.Range("K122").Value = 8692.636200000001
$ws.Range("M122").Value = -6242.636200000001
$ws.Range("H132").Value = 2394.5833
$ws.Range("I132").Value = 1526.2222
$ws.Range("K132").Value = 4578.6666
$ws.Range("M132").Value = -2048.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("H99").Value = 2140.8
$ws.Range("I99").Value = 1612.125
$ws.Range("J99").Value = 4255.5
$ws.Range("K99").Value = 1612.125
$ws.Range("L99").Value = 4255.5
$ws.Range("M99").Value = -114.125
$ws.Range("N99").Value = -7251.5
$ws.Range("H107").Value = 7308.6665
$ws.Range("I107").Value = 1100
$ws.Range("J107").Value = 10413
$ws.Range("K107").Value = 1100
$ws.Range("L107").Value = 10413
$ws.Range("M107").Value = 820
$ws.Range("N107").Value = -14253
$ws.Range("H123").Value = 100000
$ws.Range("J123").Value = 100000
$ws.Range("L123").Value = 100000
$ws.Range("N123").Value = -109800
$ws.Range("H134").Value = 2684.853
$ws.Range("I134").Value = 2264.3635
$ws.Range("K134").Value = 6793.0905
$ws.Range("M134").Value = -4258.0905
$ws.Range("H138").Value = 99998.5
$ws.Range("J138").Value = 99998.5
$ws.Range("L138").Value = 99998.5
$ws.Range("N138").Value = -110278.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12453.434
$ws.Range("I58").Value = 1595.12
$ws.Range("J58").Value = 66745
$ws.Range("K58").Value = 1595.12
$ws.Range("L58").Value = 66745
$ws.Range("M58").Value = -1392.12
$ws.Range("N58").Value = -67151
$ws.Range("H94").Value = 1941.8462
$ws.Range("I94").Value = 1818
$ws.Range("J94").Value = 1996.8889
$ws.Range("K94").Value = 1818
$ws.Range("L94").Value = 1996.8889
$ws.Range("M94").Value = -1367
$ws.Range("N94").Value = -2898.8889
$ws.Range("H105").Value = 2741.2856
$ws.Range("I105").Value = 2948.1667
$ws.Range("K105").Value = 2948.1667
$ws.Range("M105").Value = -1201.1667
$ws.Range("H136").Value = 12453.434
$ws.Range("I136").Value = 1595.12
$ws.Range("J136").Value = 66745
$ws.Range("K136").Value = 4785.36
$ws.Range("L136").Value = 200235
$ws.Range("M136").Value = -2235.36
$ws.Range("N136").Value = -205335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 126.208336
$ws.Range("I12").Value = 236.4
$ws.Range("K12").Value = 709.2
$ws.Range("M12").Value = -536.2
$ws.Range("H131").Value = 1693
$ws.Range("J131").Value = 2998.8
$ws.Range("L131").Value = 8996.400000000001
$ws.Range("N131").Value = -19076.4
$ws.Range("H139").Value = 2696
$ws.Range("I139").Value = 2366.8572
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 7100.571599999999
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = -1960.571599999999
$ws.Range("N139").Value = -25280
$ws.Range("H141").Value = 5408.8667
$ws.Range("I141").Value = 4145.4546
$ws.Range("J141").Value = 8883.25
$ws.Range("K141").Value = 12436.3638
$ws.Range("L141").Value = 26649.75
$ws.Range("M141").Value = -7256.363799999999
$ws.Range("N141").Value = -37009.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3608.4849
$ws.Range("I122").Value = 3264.6428
$ws.Range("K122").Value = 9793.928400000001
$ws.Range("M122").Value = -7343.928400000001
$ws.Range("H126").Value = 3372.25
$ws.Range("I126").Value = 1996.5
$ws.Range("K126").Value = 5989.5
$ws.Range("M126").Value = -3519.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2715.55
$ws.Range("I46").Value = 774.8333
$ws.Range("J46").Value = 5626.625
$ws.Range("K46").Value = 774.8333
$ws.Range("L46").Value = 5626.625
$ws.Range("M46").Value = -586.8333
$ws.Range("N46").Value = -6002.625
$ws.Range("H61").Value = 1350.375
$ws.Range("I61").Value = 551
$ws.Range("J61").Value = 2149.75
$ws.Range("K61").Value = 551
$ws.Range("L61").Value = 2149.75
$ws.Range("M61").Value = -349
$ws.Range("N61").Value = -2553.75
$ws.Range("H113").Value = 1350.375
$ws.Range("I113").Value = 551
$ws.Range("J113").Value = 2149.75
$ws.Range("K113").Value = 551
$ws.Range("L113").Value = 2149.75
$ws.Range("M113").Value = 1619
$ws.Range("N113").Value = -6489.75
$ws.Range("H122").Value = 3939.138
$ws.Range("I122").Value = 2907.0908
$ws.Range("J122").Value = 7182.7144
$ws.Range("K122").Value = 8721.2724
$ws.Range("L122").Value = 21548.1432
$ws.Range("M122").Value = -6271.2724
$ws.Range("N122").Value = -26448.1432
$ws.Range("H132").Value = 1743.6207
$ws.Range("I132").Value = 1379
$ws.Range("K132").Value = 4137
$ws.Range("M132").Value = -1607
$ws.Range("H136").Value = 6664.7856
$ws.Range("I136").Value = 6240
$ws.Range("K136").Value = 18720
$ws.Range("M136").Value = -16170

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1646.081
$ws.Range("I132").Value = 936.48
$ws.Range("K132").Value = 2809.44
$ws.Range("M132").Value = -279.4400000000001
$ws.Range("H137").Value = 109785
$ws.Range("J137").Value = 109785
$ws.Range("L137").Value = 109785
$ws.Range("N137").Value = -119985

